# Add a "2022-Q4" quarter sheet with fund-holding detail data, insert it
# right after the "总计" (total) summary sheet, and update the "总计"
# sheet's table with the new quarter's aggregate row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" worksheet right after "总计".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q4 = $wb.Worksheets.Add($null, $totalSheet)
$q4.Name = "2022-Q4"

# Page margins matching the other quarterly detail sheets (0.75in / 1in / 0.5in).
$q4.PageSetup.LeftMargin = 54
$q4.PageSetup.RightMargin = 54
$q4.PageSetup.TopMargin = 72
$q4.PageSetup.BottomMargin = 72
$q4.PageSetup.HeaderMargin = 36
$q4.PageSetup.FooterMargin = 36

# Header row (B1:H1) - bold, boxed border, centered/top aligned.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 2
    $cell = $q4.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# Data rows 2-6: A=index (numeric, styled), B=code, C=name, D=scale,
# E=position, F=pct, G=value, H=rank(numeric). B-G are stored as text,
# matching how the other quarterly sheets hold these numeric-looking
# values as text.
$rows = @(
    @("001322", "东吴新趋势价值线灵活配置混合", "1.64", "93.49", "8.95", "0.1468", 3),
    @("014376", "东吴新能源汽车股票A", "0.87", "93.11", "9.12", "0.0793", 4),
    @("001323", "东吴移动互联灵活配置混合A", "0.56", "93.45", "8.49", "0.0475", 3),
    @("014377", "东吴新能源汽车股票C", "0.25", "93.11", "9.12", "0.0228", 4),
    @("002170", "东吴移动互联灵活配置混合C", "0.05", "93.45", "8.49", "0.0042", 3)
)

# Force columns B:G to text so numeric-looking strings (fund codes,
# percentages) keep leading zeros / exact text form instead of being
# coerced to numbers.
$q4.Range("B2:G6").NumberFormat = "@"

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $r + 2
    $data = $rows[$r]

    $aCell = $q4.Cells.Item($rowNum, 1)
    $aCell.Value = $r
    $aCell.Font.Bold = $true
    $aCell.Borders.LineStyle = 1
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160

    $q4.Cells.Item($rowNum, 2).Value = $data[0]
    $q4.Cells.Item($rowNum, 3).Value = $data[1]
    $q4.Cells.Item($rowNum, 4).Value = $data[2]
    $q4.Cells.Item($rowNum, 5).Value = $data[3]
    $q4.Cells.Item($rowNum, 6).Value = $data[4]
    $q4.Cells.Item($rowNum, 7).Value = $data[5]
    $q4.Cells.Item($rowNum, 8).Value = $data[6]
}

# ---------------------------------------------------------------------
# 2. Rewrite the "总计" sheet's table so a 2022-Q4 row leads, and every
#    other quarter's row shifts down by one.
# ---------------------------------------------------------------------
$summaryRows = @(
    @("2022-Q4", 5, 0.3),
    @("2022-Q3", 7, 0.38),
    @("2022-Q2", 8, 0.47),
    @("2022-Q1", 13, 1.46),
    @("2021-Q4", 17, 4.39),
    @("2021-Q3", 5, 0.22),
    @("2021-Q2", 19, 3.31),
    @("2021-Q1", 12, 18.69),
    @("2020-Q4", 1, 7.73)
)

for ($r = 0; $r -lt $summaryRows.Length; $r++) {
    $rowNum = $r + 2
    $data = $summaryRows[$r]

    $aCell = $totalSheet.Cells.Item($rowNum, 1)
    $aCell.Value = $r
    $aCell.Font.Bold = $true
    $aCell.Borders.LineStyle = 1
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160

    $totalSheet.Cells.Item($rowNum, 2).Value = $data[0]
    $totalSheet.Cells.Item($rowNum, 3).Value = $data[1]
    $totalSheet.Cells.Item($rowNum, 4).Value = $data[2]
}

# ---------------------------------------------------------------------
# 3. Restore "总计" as the displayed/active sheet (adding a sheet makes
#    the new sheet active by default).
# ---------------------------------------------------------------------
$totalSheet.Activate()
